# Updated symbol list on Sat Feb  4 10:48:27 UTC 2023 with GitHub Actions
# Applies the per-cell value updates described in the target diff for sheet1
# (cryptos price/volume table). Rows 6 and 7 also swap their Coin/Link text
# (FTXToken <-> KuCoinToken) in addition to their own Price/Volume updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'329.70"
$ws.Cells.Item(2, 5).Value = "'1.40%"
$ws.Cells.Item(3, 5).Value = "'4.10%"
$ws.Cells.Item(4, 4).Value = "'5.632"
$ws.Cells.Item(4, 5).Value = "'-0.56%"
$ws.Cells.Item(5, 4).Value = "'0.08208"
$ws.Cells.Item(5, 5).Value = "'2.21%"
$ws.Cells.Item(6, 2).Value = "KuCoinToken"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(6, 4).Value = "'8.756"
$ws.Cells.Item(6, 5).Value = "'1.60%"
$ws.Cells.Item(7, 2).Value = "FTXToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(7, 4).Value = "'2.011"
$ws.Cells.Item(7, 5).Value = "'-0.35%"
$ws.Cells.Item(8, 5).Value = "'0.36%"
$ws.Cells.Item(9, 5).Value = "'1.82%"
$ws.Cells.Item(10, 4).Value = "'0.9230"
$ws.Cells.Item(10, 5).Value = "'0.07%"
$ws.Cells.Item(11, 4).Value = "'0.1278"
$ws.Cells.Item(11, 5).Value = "'2.95%"
$ws.Cells.Item(12, 4).Value = "'0.1956"
$ws.Cells.Item(12, 5).Value = "'-0.48%"
$ws.Cells.Item(13, 4).Value = "'0.09376"
$ws.Cells.Item(13, 5).Value = "'2.28%"
$ws.Cells.Item(14, 4).Value = "'0.03844"
$ws.Cells.Item(14, 5).Value = "'7.92%"
$ws.Cells.Item(15, 5).Value = "'0.96%"
$ws.Cells.Item(16, 4).Value = "'0.001306"
$ws.Cells.Item(16, 5).Value = "'0.46%"
$ws.Cells.Item(17, 4).Value = "'0.006124"
$ws.Cells.Item(17, 5).Value = "'0.54%"
$ws.Cells.Item(19, 4).Value = "'3.446"
$ws.Cells.Item(19, 5).Value = "'2.84%"
$ws.Cells.Item(21, 4).Value = "'8.343"
$ws.Cells.Item(21, 5).Value = "'-4.36%"
$ws.Cells.Item(22, 5).Value = "'-0.34%"
$ws.Cells.Item(23, 4).Value = "'0.2661"
$ws.Cells.Item(23, 5).Value = "'6.32%"
$ws.Cells.Item(24, 4).Value = "'0.04393"
$ws.Cells.Item(24, 5).Value = "'0.51%"
$ws.Cells.Item(25, 5).Value = "'-0.23%"
$ws.Cells.Item(26, 4).Value = "'0.004314"
$ws.Cells.Item(26, 5).Value = "'-6.34%"
$ws.Cells.Item(27, 4).Value = "'0.0001200"
$ws.Cells.Item(27, 5).Value = "'-2.43%"
$ws.Cells.Item(39, 4).Value = "'0.02754"
$ws.Cells.Item(39, 5).Value = "'9.35%"
$ws.Cells.Item(40, 4).Value = "'0.05476"
$ws.Cells.Item(40, 5).Value = "'2.71%"
$ws.Cells.Item(41, 4).Value = "'0.007882"
$ws.Cells.Item(41, 5).Value = "'5.30%"
$ws.Cells.Item(42, 4).Value = "'0.1421"
$ws.Cells.Item(42, 5).Value = "'1.08%"
$ws.Cells.Item(43, 4).Value = "'0.008945"
$ws.Cells.Item(43, 5).Value = "'-9.73%"
$ws.Cells.Item(44, 4).Value = "'0.002170"
$ws.Cells.Item(44, 5).Value = "'2.59%"
$ws.Cells.Item(45, 4).Value = "'0.01140"
$ws.Cells.Item(45, 5).Value = "'2.92%"
$ws.Cells.Item(46, 4).Value = "'0.00006770"
$ws.Cells.Item(46, 5).Value = "'1.38%"
$ws.Cells.Item(47, 5).Value = "'0.01%"
$ws.Cells.Item(48, 4).Value = "'0.003189"
$ws.Cells.Item(48, 5).Value = "'7.22%"
$ws.Cells.Item(49, 5).Value = "'0.07%"
$ws.Cells.Item(50, 4).Value = "'0.00002100"
$ws.Cells.Item(50, 5).Value = "'0.01%"
$ws.Cells.Item(51, 4).Value = "'0.0002000"
$ws.Cells.Item(51, 5).Value = "'0.01%"

